$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" '24.618.65'
Set-TextValue "E2" '  -1.41%  '

Set-TextValue "D3" '1.675.13'
Set-TextValue "E3" '  -2.02%  '

Set-TextValue "D4" '1.004'
Set-TextValue "E4" '  +0.30%  '

Set-TextValue "D5" '313.97'
Set-TextValue "E5" '  -1.15%  '

Set-TextValue "E6" '  +0.25%  '

Set-TextValue "D7" '0.3891'
Set-TextValue "E7" '  -3.44%  '

Set-TextValue "D8" '0.3936'
Set-TextValue "E8" '  -3.49%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D9" '1.004'
Set-TextValue "E9" '  +0.32%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D10" '51.74'
Set-TextValue "E10" '  -4.09%  '

Set-TextValue "D11" '1.390'
Set-TextValue "E11" '  -6.27%  '

Set-TextValue "D12" '0.08624'
Set-TextValue "E12" '  -2.40%  '

Set-TextValue "D13" '25.12'
Set-TextValue "E13" '  -4.88%  '

Set-TextValue "D14" '7.289'
Set-TextValue "E14" '  -3.07%  '

Set-TextValue "D15" '7.762'
Set-TextValue "E15" '  -4.74%  '

Set-TextValue "D16" '0.00001311'
Set-TextValue "E16" '  -3.79%  '

Set-TextValue "D17" '1.713.03'
Set-TextValue "E17" '  -2.27%  '

Set-TextValue "D18" '93.58'
Set-TextValue "E18" '  -3.56%  '

Set-TextValue "D19" '0.07051'
Set-TextValue "E19" '  -1.76%  '

Set-TextValue "D20" '20.59'
Set-TextValue "E20" '  -2.68%  '

Set-TextValue "D21" '7.053'
Set-TextValue "E21" '  -3.02%  '

Set-TextValue "E22" '  +0.45%  '

Set-TextValue "D23" '14.00'
Set-TextValue "E23" '  -2.93%  '

Set-TextValue "D24" '24.633.31'
Set-TextValue "E24" '  -1.35%  '

Set-TextValue "D25" '2.363'
Set-TextValue "E25" '  +1.85%  '

Set-TextValue "D26" '23.12'
Set-TextValue "E26" '  -0.85%  '

Set-TextValue "D27" '2.717'
Set-TextValue "E27" '  -6.79%  '

Set-TextValue "D28" '162.16'
Set-TextValue "E28" '  -2.94%  '

Set-TextValue "D29" '5.746'
Set-TextValue "E29" '  -8.05%  '

Set-TextValue "D30" '146.47'
Set-TextValue "E30" '  -0.30%  '

Set-TextValue "D31" '7.882'
Set-TextValue "E31" '  -6.33%  '

Set-TextValue "D32" '2.547'
Set-TextValue "E32" '  +14.31%  '

Set-TextValue "D33" '1.878.86'
Set-TextValue "E33" '  -1.70%  '

Set-TextValue "D34" '0.08342'
Set-TextValue "E34" '  -6.17%  '

Set-TextValue "D35" '0.03035'
Set-TextValue "E35" '  -6.36%  '

Set-TextValue "D36" '0.2813'
Set-TextValue "E36" '  -1.97%  '

Set-TextValue "D37" '6.836'
Set-TextValue "E37" '  -6.06%  '

Set-TextValue "D38" '0.9822'
Set-TextValue "E38" '  -4.81%  '

Set-TextValue "D39" '0.09484'
Set-TextValue "E39" '  +1.36%  '

Set-TextValue "D40" '1.544'
Set-TextValue "E40" '  +4.42%  '

Set-TextValue "D41" '10.50'
Set-TextValue "E41" '  -3.64%  '

Set-TextValue "D42" '0.7875'
Set-TextValue "E42" '  -7.23%  '

Set-TextValue "D43" '13.50'
Set-TextValue "E43" '  -5.81%  '

Set-TextValue "D44" '16.47'
Set-TextValue "E44" '  -6.67%  '

Set-TextValue "D45" '0.7104'
Set-TextValue "E45" '  -4.83%  '

Set-TextValue "D46" '2.552'
Set-TextValue "E46" '  -6.95%  '

Set-TextValue "D47" '4.186'
Set-TextValue "E47" '  -1.60%  '

Set-TextValue "D48" '0.08647'
Set-TextValue "E48" '  +3.11%  '

Set-TextValue "D49" '1.002'
Set-TextValue "E49" '  +0.26%  '

Set-TextValue "D50" '1.318'
Set-TextValue "E50" '  -5.90%  '

Set-TextValue "D51" '137.02'
Set-TextValue "E51" '  -3.74%  '

Write-Host "Applied cryptos list update"